$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 14.21340333333333
$ws.Range("H2").Value = 42.64021
$ws.Range("I2").Value = 0.07497543485230342
$ws.Range("J2").Value = 0.07497543485230343
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.09934133333334
$ws.Range("N2").Value = 63.29802400000001
$ws.Range("O2").Value = 0.2917236204149438
$ws.Range("P2").Value = 0.2917236204149438
$ws.Range("Q2").Value = 299.8934484383378
$ws.Range("R2").Value = 2699.04103594504
$ws.Range("S2").Value = 0.02187210529729871
$ws.Range("T2").Value = 0.02187210529729872

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 14.21340333333333
$ws.Range("H3").Value = 42.64021
$ws.Range("I3").Value = 0.07497543485230342
$ws.Range("J3").Value = 0.07497543485230343
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.81943766666667
$ws.Range("N3").Value = 107.458313
$ws.Range("O3").Value = 0.4952465516465762
$ws.Range("P3").Value = 0.4952465516465762
$ws.Range("Q3").Value = 509.1161147295255
$ws.Range("R3").Value = 4582.045032565729
$ws.Range("S3").Value = 0.03713132556880579
$ws.Range("T3").Value = 0.0371313255688058

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 14.21340333333333
$ws.Range("H4").Value = 42.64021
$ws.Range("I4").Value = 0.07497543485230342
$ws.Range("J4").Value = 0.07497543485230343
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 15.40769666666667
$ws.Range("N4").Value = 46.22309
$ws.Range("O4").Value = 0.2130298279384801
$ws.Range("P4").Value = 0.2130298279384801
$ws.Range("Q4").Value = 218.9958071609889
$ws.Range("R4").Value = 1970.9622644489
$ws.Range("S4").Value = 0.01597200398619892
$ws.Range("T4").Value = 0.01597200398619892

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.428335
$ws.Range("H5").Value = 88.285005
$ws.Range("I5").Value = 0.1552339127976335
$ws.Range("J5").Value = 0.1552339127976336
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.09934133333334
$ws.Range("N5").Value = 63.29802400000001
$ws.Range("O5").Value = 0.2917236204149438
$ws.Range("P5").Value = 0.2917236204149438
$ws.Range("Q5").Value = 620.9184850366801
$ws.Range("R5").Value = 5588.266365330121
$ws.Range("S5").Value = 0.04528539905250333
$ws.Range("T5").Value = 0.04528539905250335

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 29.428335
$ws.Range("H6").Value = 88.285005
$ws.Range("I6").Value = 0.1552339127976335
$ws.Range("J6").Value = 0.1552339127976336
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.81943766666667
$ws.Range("N6").Value = 107.458313
$ws.Range("O6").Value = 0.4952465516465762
$ws.Range("P6").Value = 0.4952465516465762
$ws.Range("Q6").Value = 1054.106411166285
$ws.Range("R6").Value = 9486.957700496565
$ws.Range("S6").Value = 0.07687906001163332
$ws.Range("T6").Value = 0.07687906001163335

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 29.428335
$ws.Range("H7").Value = 88.285005
$ws.Range("I7").Value = 0.1552339127976335
$ws.Range("J7").Value = 0.1552339127976336
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 15.40769666666667
$ws.Range("N7").Value = 46.22309
$ws.Range("O7").Value = 0.2130298279384801
$ws.Range("P7").Value = 0.2130298279384801
$ws.Range("Q7").Value = 453.42285908505
$ws.Range("R7").Value = 4080.80573176545
$ws.Range("S7").Value = 0.03306945373349689
$ws.Range("T7").Value = 0.03306945373349691

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 145.9323983333333
$ws.Range("H8").Value = 437.797195
$ws.Range("I8").Value = 0.769790652350063
$ws.Range("J8").Value = 0.769790652350063
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.09934133333334
$ws.Range("N8").Value = 63.29802400000001
$ws.Range("O8").Value = 0.2917236204149438
$ws.Range("P8").Value = 0.2917236204149438
$ws.Range("Q8").Value = 3079.077484026965
$ws.Range("R8").Value = 27711.69735624268
$ws.Range("S8").Value = 0.2245661160651417
$ws.Range("T8").Value = 0.2245661160651418

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 145.9323983333333
$ws.Range("H9").Value = 437.797195
$ws.Range("I9").Value = 0.769790652350063
$ws.Range("J9").Value = 0.769790652350063
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 35.81943766666667
$ws.Range("N9").Value = 107.458313
$ws.Range("O9").Value = 0.4952465516465762
$ws.Range("P9").Value = 0.4952465516465762
$ws.Range("Q9").Value = 5227.216445648004
$ws.Range("R9").Value = 47044.94801083203
$ws.Range("S9").Value = 0.3812361660661371
$ws.Range("T9").Value = 0.3812361660661371

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 145.9323983333333
$ws.Range("H10").Value = 437.797195
$ws.Range("I10").Value = 0.769790652350063
$ws.Range("J10").Value = 0.769790652350063
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.40769666666667
$ws.Range("N10").Value = 46.22309
$ws.Range("O10").Value = 0.2130298279384801
$ws.Range("P10").Value = 0.2130298279384801
$ws.Range("Q10").Value = 2248.482127359172
$ws.Range("R10").Value = 20236.33914623255
$ws.Range("S10").Value = 0.1639883702187843
$ws.Range("T10").Value = 0.1639883702187843
